$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "<<PK>> StaffID: int   increment" -- merge the separate ": " run
# and the "int   increment" run that follow the StaffID spell-check markers
# into a single run reading ": int   increment".
# ---------------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*StaffID: int*increment*") {
        $p.Range.Find.Execute(": int   increment", $false, $false, $false, $false, $false, $true, 1, $false, ": int   increment", 2) | Out-Null
    }
}

# ---------------------------------------------------------------------------
# Change 2: drop the "Will there be a table for storing media..." paragraph
# (and one of the two blank paragraphs that precede it), replacing them with
# a new "Interests" table followed by a single blank paragraph.
# ---------------------------------------------------------------------------
$lastTable = $d.Tables.Item($d.Tables.Count)
$endOfTable = $lastTable.Range.End

$needle = "Will there be a table for storing media (Videos) ?????   Do we need to ask the client about it ?"
$searchRng = $d.Range($endOfTable, $d.Content.End)
$searchRng.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$sentenceStart = $searchRng.Start
$sentenceEnd = $searchRng.End
$blankCount = $sentenceStart - $endOfTable

# Remove the sentence text itself (pure in-paragraph text, safe to delete).
$sentenceRng = $d.Range($sentenceStart, $sentenceEnd)
$sentenceRng.Delete()

# Remove every blank paragraph that sat between the table and the sentence
# paragraph -- each is a lone paragraph mark, deleted one at a time so the
# merge behaves like pressing Delete at the start of the line. (The former
# sentence paragraph, now textless, becomes the sole survivor and is
# cleaned up below.)
for ($k = 1; $k -le $blankCount; $k++) {
    $markRng = $d.Range($endOfTable, $endOfTable + 1)
    $markRng.Delete()
}

# The former sentence paragraph is now empty text but may still carry
# leftover proofing-error markup (<w:proofErr/>) with no run behind it;
# replace it outright with a clean empty paragraph.
$trailingRng = $d.Range($endOfTable, $d.Content.End)
$trailingRng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# Build the new "Interests" table as a literal OOXML fragment so it matches
# the source table formatting exactly (shaded header cell, spell-check
# markers around the field names, etc.).
$interestsTableXml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="9350"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9350" w:type="dxa"/><w:shd w:val="clear" w:color="auto" w:fill="AEAAAA" w:themeFill="background2" w:themeFillShade="BF"/></w:tcPr><w:p><w:r><w:t>Interests</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9350" w:type="dxa"/></w:tcPr><w:p><w:r><w:t xml:space="preserve">&lt;&lt;PK&gt;&gt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>InterestID</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> int    increment</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">               </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>InterestDesc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Characters {size = 250}</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">               </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>CreatedDate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: date</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">               </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LastModifiedDate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: date</w:t></w:r></w:p></w:tc></w:tr></w:tbl>
'@

# Insert a fresh blank paragraph followed by the table right before the
# blank paragraph that is left over from the sentence, giving:
#   </w:tbl><w:p/><w:tbl>...Interests...</w:tbl><w:p/><w:sectPr>
$insertionPoint = $d.Range($endOfTable, $endOfTable)
$insertionPoint.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' + $interestsTableXml)
